$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8171243547524514, 1.0904291384304117, 0.9321738139348509, 2.66060724226815)"
$ws.Range("C2").Value = "JSU(-0.9163983036770068, 1.253619990026913, 4.2036286534555884, 6.184489716725249)"
$ws.Range("D2").Value = "NIG(1.014749777416181, 0.8168464095274359, 0.5502315334891282, 3.119678845059713)"
$ws.Range("E2").Value = "NIG(1.9709875769239658, 1.605681370657908, 2.2036006119986884, 6.893476661133576)"
